$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated stats for 2025-12 (row 25)
$ws.Range("C25").Value = 1011
$ws.Range("D25").Value = 6043915
$ws.Range("E25").Value = 932.7029320987655
$ws.Range("G25").Value = 7.782515991471217
$ws.Range("H25").Value = 26.5782184388774
